# Saldo_guide.xlsx update: roll the extract date forward one day
# (2024-10-28 -> 2024-10-29) and refresh the sheet name + a handful of
# balance values that changed between the two extracts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new extraction timestamp.
$ws.Name = "IClientBalance-20241029-084130-"

# Column G holds the extract date serial for every data row (2-274):
# 45593 (2024-10-28) -> 45594 (2024-10-29).
$ws.Range("G2:G274").Value = 45594

# A handful of rows had their balance (columns E and H) restated between
# the two extracts. Update both columns together, per row.
$balanceUpdates = @{
    15  = 362.76
    52  = 194.66
    57  = 862.8
    107 = 1111.37
    109 = 6822.72
    112 = 782.38
    232 = 910.3
}

foreach ($row in $balanceUpdates.Keys) {
    $value = $balanceUpdates[$row]
    $ws.Range("E$row").Value = $value
    $ws.Range("H$row").Value = $value
}
